$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2021 column (R) of data, mirroring the formatting that column
# Q (2020) already has for each row.
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 99.4
$ws.Range("R6").Value = 98.1
$ws.Range("R7").Value = 99.319469393395053
$ws.Range("R8").Value = 99.442213297634979
$ws.Range("R9").Value = 99.1
$ws.Range("R10").Value = 99.3
$ws.Range("R11").Value = 99.799160124155549
$ws.Range("R12").Value = 99.3
$ws.Range("R13").Value = 99.538370126605429
$ws.Range("R14").Value = 99.765563948945029

# The active selection in the author's edit ended up on U4.
$ws.Range("U4").Select()
